# Applies the "overwrite old files with RMI version" update:
#  - About sheet: new "New Mexico" header cell, updated date, source citation
#    refreshed to the 2021 draft GHG inventory, and selection moved to B8.
#  - Data sheet: excerpt years rolled from 2009-2013 to 2015-2019 with new
#    CO2/CH4/N2O figures (the ratio formulas recalc automatically), and
#    selection moved to F6.
#  - RPEpUACE sheet: selection range extended from B2 to B2:B13 (values
#    there are formula-driven off Data! and recalc automatically).

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("Data")
$wsRPE = $wb.Worksheets.Item("RPEpUACE")

# ---- About sheet ----
$wsAbout.Range("B1").Value = "New Mexico"
$wsAbout.Range("C1").Value = 44515

$wsAbout.Range("B3").Value = "US EPA"
$wsAbout.Range("B4").Value = 2021
$wsAbout.Range("B5").Value = "Draft Inventory of US Greenhouse Gas Emissions Emissions and Sinks"
$wsAbout.Range("B6").Value = "https://www.epa.gov/sites/production/files/2021-02/documents/us-ghg-inventory-2021-main-text.pdf"
$wsAbout.Range("B7").Value = "Table 6-3"

# ---- Data sheet ----
$wsData.Range("A1").Value = "Excerpt from Table 6-2:"

$wsData.Range("B3").Value = 2015
$wsData.Range("C3").Value = 2016
$wsData.Range("D3").Value = 2017
$wsData.Range("E3").Value = 2018
$wsData.Range("F3").Value = 2019

$wsData.Range("B4").Value = -791695
$wsData.Range("C4").Value = -855998
$wsData.Range("D4").Value = -792046
$wsData.Range("E4").Value = -824885
$wsData.Range("F4").Value = -812695

$wsData.Range("B5").Value = 663
$wsData.Range("C5").Value = 308
$wsData.Range("D5").Value = 614
$wsData.Range("E5").Value = 552
$wsData.Range("F5").Value = 552

$wsData.Range("B6").Value = 38
$wsData.Range("C6").Value = 18
$wsData.Range("D6").Value = 36
$wsData.Range("E6").Value = 32
$wsData.Range("F6").Value = 32

$excel.CalculateFull()

# ---- Selections (order matters: last-selected sheet becomes the active tab) ----
$wsRPE.Range("B2:B13").Select()
$wsData.Range("F6").Select()
$wsAbout.Range("B8").Select()
